$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 0.001
$ws.Range("K10").Value = 483
$ws.Range("L10").Value = 0.00161
